# FindNewsInput.xlsx: "Updated to use work items"
#  - category sample value changed from "Features" to "Sport"
#  - worksheet page orientation set to portrait
#  - active selection reset to the top-left cell (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sample category value in row 2 from "Features" to "Sport"
# (both are valid options in the B2 data-validation drop-down list).
$ws.Range("B2").Value = "Sport"

# Set the sheet up for portrait printing.
$ws.PageSetup.Orientation = 1

# Clear the lingering C3 selection left over from editing, returning the
# cursor to the default top-left cell.
$ws.Activate()
$ws.Range("A1").Select()
